$wb = $excel.ActiveWorkbook

# --- 1. Fix up the matching row in "item_inventory" (A420 was a placeholder 0) ---
$invWs = $wb.Worksheets.Item("item_inventory")
$invWs.Activate()
$invWs.Range("A420").Value = 1901
$invWs.Rows.Item(420).Select()

# --- 2. Insert the new item row into the "item" sheet ---
$itemWs = $wb.Worksheets.Item("item")
$itemWs.Activate()
$itemWs.Rows.Item(582).Insert()

$itemWs.Range("A582").Value = 1901
$itemWs.Range("B582").Value = "隔离缎带"
$itemWs.Range("C582").Value = "Spurning Ribbon"

$itemWs.Range("H579").Select()
